$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7-11 (they are removed entirely in the target)
$ws.Range("A7:B11").EntireRow.Delete()

# Update remaining data rows 2-6 with new values
$values = @(
    @(1, 41),
    @(4, 36),
    @(2, 35),
    @(3, 26),
    @(0, 26)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
